$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harmony_type values from "ATR" to "non-harmonic" for specific rows
$ws.Range("B4").Value = "non-harmonic"
$ws.Range("B9").Value = "non-harmonic"
$ws.Range("B24").Value = "non-harmonic"
$ws.Range("B41").Value = "non-harmonic"
$ws.Range("B51").Value = "non-harmonic"
$ws.Range("B59").Value = "non-harmonic"

# Apply AutoFilter to A1:B60 (also creates the _FilterDatabase defined name)
$ws.Range("A1:B60").AutoFilter() | Out-Null

# Make sure the sheet-scoped hidden _FilterDatabase defined name matches exactly
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$60")
$filterName.Visible = $false

# Update selection / view (also clears the previous topLeftCell scroll position)
$ws.Range("C10").Select() | Out-Null
